$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.309980630874634
$ws.Range("B1").Value = 2.411827802658081
$ws.Range("C1").Value = 6.041519641876221
$ws.Range("D1").Value = 1.667244076728821
$ws.Range("E1").Value = 1.302856922149658
